$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns (HRA, CAR) before the existing "Spcl Allowance" column (col N / 14) ---
$ws.Columns.Item(14).EntireColumn.Insert()
$ws.Columns.Item(14).EntireColumn.Insert()

# --- New header cells for the inserted columns ---
$ws.Range("N1").Value = "HRA"
$ws.Range("O1").Value = "CAR"

# --- New data cells (row 2) for the inserted columns ---
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0

# --- New data cells (row 3) for the inserted columns ---
$ws.Range("N3").Value = 800.0
$ws.Range("O3").Value = 400.0

# --- Fix row 3 designation text ---
$ws.Range("F3").Value = "Junior Development"

# --- Fix row 3 values that had erroneous extra digits ---
$ws.Range("H3").Value = 120000.0
$ws.Range("I3").Value = 10000.0
$ws.Range("J3").Value = 9345.0
$ws.Range("M3").Value = 4000.0

# --- Fix the (now shifted) row 3 values after the insertion point ---
$ws.Range("P3").Value = 4145.0
$ws.Range("Q3").Value = 0.0
$ws.Range("R3").Value = 9345.0
$ws.Range("S3").Value = 480.0
$ws.Range("T3").Value = 163.54
$ws.Range("U3").Value = 0.0
$ws.Range("V3").Value = 0.0
$ws.Range("W3").Value = 643.54
$ws.Range("X3").Value = 8701.46

# --- Column widths (closest achievable approximation of the target best-fit widths;
#     the host's pixel-grid rounds ColumnWidth to the nearest 1/7 character, so these
#     inputs are chosen to land on the nearest representable value to the target) ---
$ws.Columns.Item(6).ColumnWidth = 16.714285714285715
$ws.Columns.Item(7).ColumnWidth = 13.428571428571429
$ws.Columns.Item(8).ColumnWidth = 10.142857142857142
$ws.Columns.Item(9).ColumnWidth = 20.0
$ws.Columns.Item(10).ColumnWidth = 13.428571428571429
$ws.Columns.Item(13).ColumnWidth = 9.0
$ws.Columns.Item(14).ColumnWidth = 6.714285714285714
$ws.Columns.Item(15).ColumnWidth = 6.714285714285714
$ws.Columns.Item(16).ColumnWidth = 11.142857142857142
$ws.Columns.Item(17).ColumnWidth = 4.571428571428571
$ws.Columns.Item(18).ColumnWidth = 9.0
$ws.Columns.Item(19).ColumnWidth = 6.714285714285714
$ws.Columns.Item(20).ColumnWidth = 7.857142857142857
$ws.Columns.Item(21).ColumnWidth = 4.571428571428571
$ws.Columns.Item(22).ColumnWidth = 5.714285714285714
$ws.Columns.Item(23).ColumnWidth = 11.142857142857142
$ws.Columns.Item(24).ColumnWidth = 9.0
